$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("Q2").Value = 2.4
$ws.Range("R2").Value = 1.53
$ws.Range("H3").Value = 3.6
$ws.Range("I3").Value = 5.25
$ws.Range("J3").Value = 2.4
$ws.Range("S3").Value = 1.44
$ws.Range("T3").Value = 2.63
$ws.Range("U3").Value = 2
$ws.Range("V3").Value = 1.75
$ws.Range("X3").Value = 7.5
$ws.Range("Y3").Value = 8.5
$ws.Range("AE3").Value = 17
$ws.Range("AI3").Value = 17
$ws.Range("AM3").Value = 401
$ws.Range("AN3").Value = 3.6
$ws.Range("AO3").Value = 9.5
$ws.Range("AT3").Value = 2.63
$ws.Range("BC3").Value = 151
$ws.Range("G4").Value = 2
$ws.Range("I4").Value = 4.75
$ws.Range("K4").Value = 1.83
$ws.Range("L4").Value = 5.5
$ws.Range("M4").Value = 1.17
$ws.Range("N4").Value = 5
$ws.Range("S4").Value = 1.67
$ws.Range("T4").Value = 2.1
$ws.Range("U4").Value = 2.5
$ws.Range("V4").Value = 1.5
$ws.Range("W4").Value = 4.75
$ws.Range("X4").Value = 7.5
$ws.Range("Y4").Value = 10
$ws.Range("AB4").Value = 41
$ws.Range("AC4").Value = 5
$ws.Range("AD4").Value = 6
$ws.Range("AE4").Value = 23
$ws.Range("AF4").Value = 101
$ws.Range("AH4").Value = 21
$ws.Range("AR4").Value = 81
$ws.Range("AS4").Value = 401
$ws.Range("AT4").Value = 2.1
$ws.Range("AU4").Value = 11
$ws.Range("AX4").Value = 29
$ws.Range("AY4").Value = 41
$ws.Range("AZ4").Value = 126
$ws.Range("BA4").Value = 201
$ws.Range("G8").Value = 1.44
$ws.Range("N8").Value = 12
$ws.Range("Q8").Value = 1.8
$ws.Range("R8").Value = 2
$ws.Range("AC8").Value = 12
$ws.Range("AE8").Value = 19
$ws.Range("AV8").Value = 51
$ws.Range("AX8").Value = 34
$ws.Range("BB8").Value = 301
$ws.Range("M11").Value = 1.07
$ws.Range("N11").Value = 9
